# Update "bitbucket-to-azure" sheet with new Bitbucket workspace/repo data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New workspace id (column B) applies to every data row.
$workspace = "anilgoudasb06"

# Repos (column C) to import, one per row.
$repos = @("almatasks", "app-n-pak", "casa-build-utils", "casaplotserver", "casashell")

# Target Azure namespace (column D) applies to every data row.
$namespace = "anilbharamagoudar06/Anil-devops"

for ($i = 0; $i -lt $repos.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $i + 1
    $ws.Cells.Item($row, 2).Value2 = $workspace
    $ws.Cells.Item($row, 3).Value2 = $repos[$i]
}

# Fill the namespace column last so the shared string for it is appended
# to the end of the shared-strings table (matching workbook save order).
for ($i = 0; $i -lt $repos.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value2 = $namespace
}

# Column D now needs to fit the longer namespace string (best-fit width).
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(4).ColumnWidth = 31.6

# Reflect the final selection recorded in the workbook.
$ws.Range("J8").Select() | Out-Null
